$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one new row (row 8) to the trip log, mirroring the existing rows.
# A8 is an empty-text cell (same as A2-A5/A7): a lone "'" yields a
# quote-prefixed empty string via COM, then resetting the style to Normal
# drops the quote-prefix formatting while keeping the empty text value.
$ws.Range("A8").Value = "'"
$ws.Range("A8").Style = "Normal"

$ws.Range("B8").Value = "حسن "

# C8 holds "23" as text (like C2..C7), not a number - same apostrophe trick.
$ws.Range("C8").Value = "'23"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = "ايتا"
$ws.Range("E8").Value = "الرحلة 3"
$ws.Range("F8").Value = "C5"
$ws.Range("G8").Value = "UNDP"
$ws.Range("H8").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٤٩:٤٦ م"
